$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing text formatting of column D (prices are stored as
# text, some of which look numeric, e.g. "572.61") while updating values.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "61.451.74"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "3.360.04"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "572.61"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").Value = "136.34"
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.358.90"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D10").Value = "7.47"
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("E11").Value = "  -2.13%  "
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("D13").Value = "3.935.20"
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("E15").Value = "  -3.39%  "
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("D17").Value = "3.361.90"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").Value = "61.530.87"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "13.97"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").Value = "5.84"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").Value = "9.32"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "379.08"
$ws.Range("E22").Value = "  -2.83%  "
$ws.Range("E23").Value = "  -3.66%  "
$ws.Range("D24").Value = "3.500.50"
$ws.Range("E24").Value = "  -0.95%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("D27").Value = "71.20"
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("D28").Value = "1.78"
$ws.Range("E28").Value = "  +9.68%  "
$ws.Range("D29").Value = "7.52"
$ws.Range("E29").Value = "  -2.89%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("E31").Value = "  +2.56%  "
$ws.Range("E32").Value = "  -1.89%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "23.50"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").Value = "5.21"
$ws.Range("E36").Value = "  -5.62%  "
$ws.Range("D38").Value = "165.63"
$ws.Range("E38").Value = "  +2.35%  "
$ws.Range("E39").Value = "  -2.06%  "
$ws.Range("D40").Value = "0.0761"
$ws.Range("E40").Value = "  -4.11%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("D43").Value = "0.768"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("E44").Value = "  -1.97%  "
$ws.Range("D45").Value = "41.47"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("E46").Value = "  -1.94%  "
$ws.Range("D47").Value = "24.06"
$ws.Range("E47").Value = "  -1.16%  "
$ws.Range("D48").Value = "6.82"
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("D49").Value = "23.09"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").Value = "2.373.35"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("E51").Value = "  -2.44%  "

# Restore the original (default) cell style now that the values are locked
# in as text, so we do not leave a stray number-format style behind.
$dRange.Style = "Normal"
